# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (default blue Office palette),
#                            used by the Notes Master.
#   ppt/theme/theme2.xml -> "Integral" (green/olive palette), used by the
#                            Slide Master (i.e. what the slides actually show).
#
# The authored change swaps the two themes' contents wholesale, so the
# deck's visible design (Slide Master) switches from "Integral" to the
# default "Office Theme" colors, while the Notes Master gets "Integral".
#
# PowerPoint's object model edits theme colors through
# Theme.ThemeColorScheme (12 slots: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink), addressed the same way the COM automation surface does it -
# by driving each ThemeColor.RGB. We push the "Office Theme" palette onto
# the presentation's (Slide Master) color scheme to reproduce the swap.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Target palette = the stock "Office Theme" colors (formerly theme1.xml),
# now applied where "Integral" (theme2.xml) used to be.
# Index : scheme slot : hex -> packed BGR long (R + G*256 + B*65536),
# PowerPoint's native RGB-color-long encoding.
$colorScheme.Colors(1).RGB  = 0          # dk1      000000
$colorScheme.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$colorScheme.Colors(3).RGB  = 6968388    # dk2      44546A
$colorScheme.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$colorScheme.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$colorScheme.Colors(6).RGB  = 3243501    # accent2  ED7D31
$colorScheme.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$colorScheme.Colors(8).RGB  = 49407      # accent4  FFC000
$colorScheme.Colors(9).RGB  = 12874308   # accent5  4472C4
$colorScheme.Colors(10).RGB = 4697456    # accent6  70AD47
$colorScheme.Colors(11).RGB = 12673797   # hlink    0563C1
$colorScheme.Colors(12).RGB = 7491477    # folHlink 954F72

# Note: the theme/clrScheme XML "name" attributes ("Integral" / "Office")
# aren't exposed as a settable property anywhere in the PowerPoint object
# model (ThemeColorScheme.Name and Design.Name map to other things, e.g.
# the Slide Master's display name - not the theme part's @name), so only
# the color values themselves are updated here, matching what COM
# automation can legitimately reach.
